$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.135.27"
$ws.Range("E2").Value = "'  -0.99%  "
$ws.Range("D3").Value = "'1.834.25"
$ws.Range("E3").Value = "'  -1.03%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("D5").Value = "'240.59"
$ws.Range("D6").Value = "'0.6643"
$ws.Range("E6").Value = "'  -4.11%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("D8").Value = "'0.2943"
$ws.Range("E8").Value = "'  -3.79%  "
$ws.Range("D9").Value = "'0.07364"
$ws.Range("E9").Value = "'  -3.80%  "
$ws.Range("D10").Value = "'22.75"
$ws.Range("E10").Value = "'  -3.15%  "
$ws.Range("E11").Value = "'  -0.66%  "
$ws.Range("D12").Value = "'1.821.34"
$ws.Range("E12").Value = "'  -1.74%  "
$ws.Range("D13").Value = "'5.012"
$ws.Range("E13").Value = "'  -2.49%  "
$ws.Range("D14").Value = "'0.6735"
$ws.Range("E14").Value = "'  -2.85%  "
$ws.Range("D15").Value = "'86.04"
$ws.Range("E15").Value = "'  -5.30%  "
$ws.Range("D16").Value = "'6.196"
$ws.Range("E16").Value = "'  -1.96%  "
$ws.Range("B17").Value = "'WrappedBTC"
$ws.Range("C17").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'28.908.02"
$ws.Range("E17").Value = "'  -1.78%  "
$ws.Range("B18").Value = "'ShibaInu"
$ws.Range("C18").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000008223"
$ws.Range("E18").Value = "'  -0.55%  "
$ws.Range("D19").Value = "'227.99"
$ws.Range("E19").Value = "'  -3.56%  "
$ws.Range("D20").Value = "'12.50"
$ws.Range("E20").Value = "'  -1.67%  "
$ws.Range("D21").Value = "'0.9981"
$ws.Range("E21").Value = "'  -0.18%  "
$ws.Range("D22").Value = "'7.248"
$ws.Range("E22").Value = "'  -5.39%  "
$ws.Range("D23").Value = "'0.9996"
$ws.Range("E23").Value = "'  -0.03%  "
$ws.Range("D24").Value = "'160.29"
$ws.Range("E24").Value = "'  +0.20%  "
$ws.Range("D25").Value = "'8.693"
$ws.Range("E25").Value = "'  -2.93%  "
$ws.Range("D26").Value = "'0.1394"
$ws.Range("E26").Value = "'  -5.53%  "
$ws.Range("D27").Value = "'17.99"
$ws.Range("E27").Value = "'  -1.21%  "
$ws.Range("D28").Value = "'1.502"
$ws.Range("E28").Value = "'  -1.97%  "
$ws.Range("D29").Value = "'4.197"
$ws.Range("E29").Value = "'  -1.21%  "
$ws.Range("D30").Value = "'4.074"
$ws.Range("E30").Value = "'  -1.54%  "
$ws.Range("D31").Value = "'1.187"
$ws.Range("E31").Value = "'  -1.40%  "
$ws.Range("D32").Value = "'0.05326"
$ws.Range("E32").Value = "'  +2.28%  "
$ws.Range("B33").Value = "'ImmutableX"
$ws.Range("C33").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7494"
$ws.Range("E33").Value = "'  -3.09%  "
$ws.Range("B34").Value = "'LidoDAOToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.855"
$ws.Range("E34").Value = "'  -0.81%  "
$ws.Range("E35").Value = "'  -1.18%  "
$ws.Range("D36").Value = "'2.679"
$ws.Range("E36").Value = "'  -0.39%  "
$ws.Range("D37").Value = "'1.317.98"
$ws.Range("E37").Value = "'  -0.43%  "
$ws.Range("D38").Value = "'0.01803"
$ws.Range("E38").Value = "'  -3.37%  "
$ws.Range("D39").Value = "'2.718"
$ws.Range("E39").Value = "'  -0.01%  "
$ws.Range("D40").Value = "'0.9202"
$ws.Range("E40").Value = "'  -2.20%  "
$ws.Range("D41").Value = "'5.971"
$ws.Range("E41").Value = "'  +3.19%  "
$ws.Range("D42").Value = "'0.9977"
$ws.Range("E42").Value = "'  -0.18%  "
$ws.Range("D43").Value = "'103.43"
$ws.Range("E43").Value = "'  -2.46%  "
$ws.Range("D44").Value = "'0.08092"
$ws.Range("E44").Value = "'  +16.51%  "
$ws.Range("D45").Value = "'0.00000000125"
$ws.Range("E45").Value = "'  +2.77%  "
$ws.Range("D46").Value = "'0.5167"
$ws.Range("E46").Value = "'  -1.11%  "
$ws.Range("D47").Value = "'1.956.35"
$ws.Range("E47").Value = "'  -2.11%  "
$ws.Range("D48").Value = "'63.75"
$ws.Range("E48").Value = "'  +1.26%  "
$ws.Range("D49").Value = "'1.754"
$ws.Range("E49").Value = "'  -1.47%  "
$ws.Range("D50").Value = "'9.267"
$ws.Range("E50").Value = "'  -5.05%  "
$ws.Range("D51").Value = "'0.05930"
$ws.Range("E51").Value = "'  -0.40%  "
